# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gets three new trailing columns:
#   H = date              (2012-04-25, the filing date)
#   I = legislator_name   (柯建銘)
#   J = legislator_id     (629)
#
# This mirrors the same three columns being appended to every sheet's
# dataframe in the source project; this workbook only shows the diff on
# the 股票 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票

$legislatorName = "柯建銘"
$legislatorId = 629
$filingDate = "2012-04-25"

# ---- header row (row 1): copy the look of the existing header cell (G1)
# onto the three new header cells, then fill in the labels. ----
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# ---- body rows (2-6): copy the look of an existing body cell (G col)
# onto the new cells in each row, then fill in the values. ----
$lastRow = 6
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r`:J$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    # Force the date column to be stored as literal text (matching the
    # source dataframe export) instead of Excel's automatic date-serial
    # inference for strings that look like dates.
    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value = $filingDate

    $ws.Range("I$r").Value = $legislatorName
    $ws.Range("J$r").Value = $legislatorId
}
